$wb = $excel.ActiveWorkbook

# ---------- Sheet: Summary ----------
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.6685393258426966
$ws1.Range("C2").Value = 0.6108374384236454
$ws1.Range("D2").Value = 0.9288389513108615
$ws1.Range("E2").Value = 0.736998514115899
$ws1.Range("F2").Value = 0.841248303934871
$ws1.Range("G2").Value = 0.9106058466318316
$ws1.Range("H2").Value = 0.7732609519000125
$ws1.Range("I2").Value = 496
$ws1.Range("J2").Value = 316
$ws1.Range("K2").Value = 218
$ws1.Range("L2").Value = 38

# ---------- Sheet: Classification Report ----------
$ws2 = $wb.Worksheets.Item("Classification Report")

$ws2.Range("B2").Value = 0.8515625
$ws2.Range("C2").Value = 0.4082397003745318
$ws2.Range("D2").Value = 0.5518987341772152

$ws2.Range("B3").Value = 0.6108374384236454
$ws2.Range("C3").Value = 0.9288389513108615
$ws2.Range("D3").Value = 0.736998514115899

$ws2.Range("B4").Value = 0.6685393258426966
$ws2.Range("C4").Value = 0.6685393258426966
$ws2.Range("D4").Value = 0.6685393258426966
$ws2.Range("E4").Value = 0.6685393258426966

$ws2.Range("B5").Value = 0.7311999692118227
$ws2.Range("C5").Value = 0.6685393258426966
$ws2.Range("D5").Value = 0.6444486241465571

$ws2.Range("B6").Value = 0.7311999692118227
$ws2.Range("C6").Value = 0.6685393258426966
$ws2.Range("D6").Value = 0.6444486241465571

# ---------- Sheet: Confusion Matrix ----------
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

$ws3.Range("B2").Value = 218
$ws3.Range("C2").Value = 316

$ws3.Range("B3").Value = 38
$ws3.Range("C3").Value = 496
